$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999608411128,
    0.99444660876183688,
    0.97726965831376089,
    0.97044805119437738,
    0.96063008928515858,
    0.93695853011186436,
    0.93035723854111296,
    0.92122144880646983,
    0.90946833073312039,
    0.89863459029151382,
    0.89703782890703043,
    0.8942569801676703,
    0.88296641023073985,
    0.8787982418505299,
    0.87620633306453199,
    0.87369954698078744,
    0.86999130783444545,
    0.86888236587915402,
    0.99221800764778845,
    0.97462604514355733,
    0.97009410087151737,
    0.96312490062735379,
    0.97023430027468538,
    0.95721285936289902,
    0.95075576366302905,
    0.93562496965976516,
    0.93401789361867804,
    0.92719528699917575,
    0.92284734904659316,
    0.92188848061741213,
    0.92809270848009739,
    0.9311959805570269,
    0.93284502281211656
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
